# Ghrh-Ghrhr.xlsx was regenerated with updated TPM input data.
# The old row 2 ("ECs" sending -> "ECs" target, i.e. self-signaling) is gone,
# and the row that used to be row 3 ("MuSCs" -> "ECs") is now the sole data
# row (row 2), with its derived-specificity / weight columns recalculated
# against the smaller (now single-row) group.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the old "ECs -> ECs" row; this shifts the old row 3 ("MuSCs -> ECs")
# up to become row 2.
$ws.Rows(2).Delete()

# Re-stamp the recalculated values for the remaining row (now a singleton
# group, so several specificity ratios collapse to 1).
$ws.Range("G2").Value2 = 0.09527133333333333
$ws.Range("H2").Value2 = 0.285814
$ws.Range("I2").Value2 = 1
$ws.Range("J2").Value2 = 1
$ws.Range("M2").Value2 = 0.3991683333333333
$ws.Range("N2").Value2 = 1.197505
$ws.Range("Q2").Value2 = 0.03802929934111111
$ws.Range("R2").Value2 = 0.34226369407
$ws.Range("S2").Value2 = 1
$ws.Range("T2").Value2 = 1
